# Auto-generated edit script: updates Leve market-price columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 5707
$ws.Range("I80").Value = 399.75
$ws.Range("J80").Value = 7636.909
$ws.Range("K80").Value = 1199.25
$ws.Range("L80").Value = 22910.727
$ws.Range("M80").Value = -201.25
$ws.Range("N80").Value = -24906.727

# Row 82
$ws.Range("H82").Value = 4135.875
$ws.Range("I82").Value = 452.5
$ws.Range("K82").Value = 1357.5
$ws.Range("M82").Value = -951.5

# Row 83
$ws.Range("H83").Value = 5707
$ws.Range("I83").Value = 399.75
$ws.Range("J83").Value = 7636.909
$ws.Range("K83").Value = 3597.75
$ws.Range("L83").Value = 68732.181
$ws.Range("M83").Value = 1394.25
$ws.Range("N83").Value = -78716.181

# Row 85
$ws.Range("H85").Value = 4135.875
$ws.Range("I85").Value = 452.5
$ws.Range("K85").Value = 1357.5
$ws.Range("M85").Value = 46.5

# Row 98
$ws.Range("H98").Value = 2963.9048
$ws.Range("I98").Value = 2662.1
$ws.Range("K98").Value = 2662.1
$ws.Range("M98").Value = -1164.1

# Row 106
$ws.Range("H106").Value = 1848.4
$ws.Range("I106").Value = 1026.8
$ws.Range("J106").Value = 2670
$ws.Range("K106").Value = 1026.8
$ws.Range("L106").Value = 2670
$ws.Range("M106").Value = -395.8
$ws.Range("N106").Value = -3932

# Row 122
$ws.Range("H122").Value = 2963.9048
$ws.Range("I122").Value = 2662.1
$ws.Range("K122").Value = 7986.299999999999
$ws.Range("M122").Value = -5536.299999999999

# Row 137
$ws.Range("H137").Value = 3733.318
$ws.Range("I137").Value = 3990.5483
$ws.Range("J137").Value = 3119.923
$ws.Range("K137").Value = 11971.6449
$ws.Range("L137").Value = 9359.769
$ws.Range("M137").Value = -9421.644899999999
$ws.Range("N137").Value = -14459.769

# Row 141
$ws.Range("H141").Value = 697060.8
$ws.Range("I141").Value = 1970.7142
$ws.Range("J141").Value = 1237686.5
$ws.Range("K141").Value = 5912.142599999999
$ws.Range("L141").Value = 3713059.5
$ws.Range("M141").Value = -732.1425999999992
$ws.Range("N141").Value = -3723419.5

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6988.4624
$ws.Range("I32").Value = 6169.6816
$ws.Range("J32").Value = 21399
$ws.Range("K32").Value = 6169.6816
$ws.Range("L32").Value = 21399
$ws.Range("M32").Value = -5882.6816
$ws.Range("N32").Value = -21973

# Row 61
$ws.Range("H61").Value = 2718.4075
$ws.Range("I61").Value = 2263.2144
$ws.Range("J61").Value = 3208.6155
$ws.Range("K61").Value = 2263.2144
$ws.Range("L61").Value = 3208.6155
$ws.Range("M61").Value = -2051.2144
$ws.Range("N61").Value = -3632.6155

# Row 122
$ws.Range("H122").Value = 1862.0364
$ws.Range("I122").Value = 1570.1025
$ws.Range("K122").Value = 4710.3075
$ws.Range("M122").Value = -2260.3075

# Row 136
$ws.Range("H136").Value = 2718.4075
$ws.Range("I136").Value = 2263.2144
$ws.Range("J136").Value = 3208.6155
$ws.Range("K136").Value = 6789.6432
$ws.Range("L136").Value = 9625.8465
$ws.Range("M136").Value = -4239.6432
$ws.Range("N136").Value = -14725.8465

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2387.4285
$ws.Range("I99").Value = 2362.1177
$ws.Range("K99").Value = 2362.1177
$ws.Range("M99").Value = -864.1176999999998

# Row 105
$ws.Range("H105").Value = 1486.5714
$ws.Range("I105").Value = 1389
$ws.Range("J105").Value = 1730.5
$ws.Range("K105").Value = 1389
$ws.Range("L105").Value = 1730.5
$ws.Range("M105").Value = 358
$ws.Range("N105").Value = -5224.5

# Row 134
$ws.Range("H134").Value = 2556.3965
$ws.Range("I134").Value = 2572.9285
$ws.Range("K134").Value = 7718.7855
$ws.Range("M134").Value = -5183.7855

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()

# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()

# Row 132
$ws.Range("H132").Value = 2361.2646
$ws.Range("I132").Value = 2220.2632
$ws.Range("J132").Value = 2539.8667
$ws.Range("K132").Value = 6660.7896
$ws.Range("L132").Value = 7619.6001
$ws.Range("M132").Value = -4130.7896
$ws.Range("N132").Value = -12679.6001

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 107.375
$ws.Range("I23").Value = 90.333336
$ws.Range("J23").Value = 111.30769
$ws.Range("K23").Value = 271.000008
$ws.Range("L23").Value = 333.92307
$ws.Range("M23").Value = -36.00000799999998
$ws.Range("N23").Value = -803.9230700000001

# Row 87
$ws.Range("H87").Value = 6996
$ws.Range("I87").Value = 6996
$ws.Range("K87").Value = 20988
$ws.Range("M87").Value = -19740

# Row 90
$ws.Range("H90").Value = 6996
$ws.Range("I90").Value = 6996
$ws.Range("K90").Value = 62964
$ws.Range("M90").Value = -56724

# Row 140
$ws.Range("H140").Value = 11115567
$ws.Range("I140").Value = 27778668
$ws.Range("J140").Value = 6833.1113
$ws.Range("K140").Value = 83336004
$ws.Range("L140").Value = 20499.3339
$ws.Range("M140").Value = -83330824
$ws.Range("N140").Value = -30859.3339

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 80033
$ws.Range("J52").Value = 80033
$ws.Range("L52").Value = 80033
$ws.Range("N52").Value = -80551

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 55556536
$ws.Range("I16").Value = 83334510
$ws.Range("J16").Value = 588.1667
$ws.Range("K16").Value = 83334510
$ws.Range("L16").Value = 588.1667
$ws.Range("M16").Value = -83334340
$ws.Range("N16").Value = -928.1667

# Row 46
$ws.Range("H46").Value = 1294.7106
$ws.Range("I46").Value = 962.8570999999999
$ws.Range("J46").Value = 5166.3335
$ws.Range("K46").Value = 962.8570999999999
$ws.Range("L46").Value = 5166.3335
$ws.Range("M46").Value = -774.8570999999999
$ws.Range("N46").Value = -5542.3335

# Row 68
$ws.Range("H68").Value = 3055.5557
$ws.Range("I68").Value = 1366.6666
$ws.Range("J68").Value = 6433.3335
$ws.Range("K68").Value = 1366.6666
$ws.Range("L68").Value = 6433.3335
$ws.Range("M68").Value = -617.6666
$ws.Range("N68").Value = -7931.3335

# Row 71
$ws.Range("H71").Value = 3055.5557
$ws.Range("I71").Value = 1366.6666
$ws.Range("J71").Value = 6433.3335
$ws.Range("K71").Value = 6833.333000000001
$ws.Range("L71").Value = 32166.6675
$ws.Range("M71").Value = -3089.333000000001
$ws.Range("N71").Value = -39654.6675

# Row 136
$ws.Range("H136").Value = 2351.4783
$ws.Range("I136").Value = 2242.75
$ws.Range("K136").Value = 6728.25
$ws.Range("M136").Value = -4178.25

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2663.8462
$ws.Range("I122").Value = 1693.6364
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 5080.9092
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -2630.9092
$ws.Range("N122").Value = -28900

# Row 126
$ws.Range("H126").Value = 32371.727
$ws.Range("I126").Value = 39925.652
$ws.Range("K126").Value = 119776.956
$ws.Range("M126").Value = -117306.956

# Row 132
$ws.Range("H132").Value = 16908.277
$ws.Range("I132").Value = 2913.0435
$ws.Range("J132").Value = 41669.08
$ws.Range("K132").Value = 8739.130500000001
$ws.Range("L132").Value = 125007.24
$ws.Range("M132").Value = -6209.130500000001
$ws.Range("N132").Value = -130067.24

# Row 136
$ws.Range("H136").Value = 3305.9524
$ws.Range("I136").Value = 2925.3076
$ws.Range("J136").Value = 3924.5
$ws.Range("K136").Value = 8775.9228
$ws.Range("L136").Value = 11773.5
$ws.Range("M136").Value = -6225.9228
$ws.Range("N136").Value = -16873.5
